$wb = $excel.ActiveWorkbook
$loginSheet = $wb.Worksheets.Item("Login")

# Add a new "Product" worksheet right after "Login" (becomes sheetId 2 / rId3,
# and becomes the active sheet -> workbook activeTab moves to index 1).
$productSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $loginSheet)
$productSheet.Name = "Product"

# Populate the test data: a numeric product id and the product name (new
# shared string "KIVIK").
$productSheet.Range("A1").Value = 99011429
$productSheet.Range("B1").Value = "KIVIK"

# Reuse the existing "demo123" cell's formatting (same font/number-format xf
# as cellXfs index 2) as the base style for the new row, then turn on
# word-wrap, producing a new cellXfs entry (index 3) with wrapText="true".
$loginSheet.Range("B1").Copy()
$productSheet.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$productSheet.Range("A1:B1").WrapText = $true

# Match the recorded selection on the new sheet (mirrors Login's B1 selection).
$productSheet.Range("B1").Select()
